$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three changed data values in row 6
$ws.Range("E6").Value = 6
$ws.Range("G6").Value = -3
$ws.Range("H6").Value = 13

# Move/record the active selection to E6 (reflected as <selection activeCell="E6" sqref="E6"/>)
$ws.Range("E6").Select()
